$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-23"

# Update the column header label (I1) for the "2022 (through ..)" column
$ws.Range("I1").Value = "2022 (through 06-23)"

# Update June total (row 7) for the 2022 column
$ws.Range("I7").Value = 108

# Update grand Total row (row 14) for the 2022 column
$ws.Range("I14").Value = 771
